$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GeomToMesh Attributes")
$ws.Activate()

# --- New attribute row: PW:DomainBlunt -------------------------------------------------
# Row 27 was blank (placeholder row); fill it in with the new attribute definition
# introduced for baffle/periodic-export support.
$ws.Range("A27").Value = "PW:DomainBlunt"
$ws.Range("B27").Value = "`$true or `$false"
$ws.Range("C27").Value = "Face"
$ws.Range("D27").Value = "Flag the domain as blunt for special dimension handling"

# --- Spacer row below grew taller ------------------------------------------------------
$ws.Rows.Item(29).RowHeight = 22

# --- View / selection state, matching where the author left the cursor -----------------
$win = $excel.ActiveWindow
try { $win.ScrollRow = 13 } catch {}
try { $win.ScrollColumn = 1 } catch {}
try { $win.Left = 22840 } catch {}
try { $win.Top = 7100 } catch {}

$ws.Range("D27").Select()
